$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Cells.Item(1, 1).Value = "Datos actualizados a 7 de Mayo de 2020 a las 20:04"

$ws.Cells.Item(4, 2).Value = 1277606
$ws.Cells.Item(4, 3).Value = 14514
$ws.Cells.Item(4, 4).Value = 214276
$ws.Cells.Item(4, 5).Value = 987278
$ws.Cells.Item(4, 6).Value = 15705
$ws.Cells.Item(4, 7).Value = 1253
$ws.Cells.Item(4, 8).Value = 76052

$ws.Cells.Item(12, 2).Value = 127655
$ws.Cells.Item(12, 3).Value = 1044
$ws.Cells.Item(12, 5).Value = 67676
$ws.Cells.Item(12, 7).Value = 21
$ws.Cells.Item(12, 8).Value = 8609

$ws.Cells.Item(15, 2).Value = 64817
$ws.Cells.Item(15, 3).Value = 1321
$ws.Cells.Item(15, 4).Value = 28954
$ws.Cells.Item(15, 5).Value = 31459
$ws.Cells.Item(15, 7).Value = 172
$ws.Cells.Item(15, 8).Value = 4404

$ws.Cells.Item(17, 2).Value = 54539
$ws.Cells.Item(17, 3).Value = 1552
$ws.Cells.Item(17, 4).Value = 16048
$ws.Cells.Item(17, 5).Value = 36654
$ws.Cells.Item(17, 7).Value = 52
$ws.Cells.Item(17, 8).Value = 1837

$ws.Cells.Item(21, 1).Value = "Ecuador"
$ws.Cells.Item(21, 2).Value = 30298
$ws.Cells.Item(21, 3).Value = 878
$ws.Cells.Item(21, 4).Value = 3433
$ws.Cells.Item(21, 5).Value = 25211
$ws.Cells.Item(21, 6).Value = 156
$ws.Cells.Item(21, 7).Value = 36
$ws.Cells.Item(21, 8).Value = 1654

$ws.Cells.Item(22, 1).Value = "Suiza"
$ws.Cells.Item(22, 2).Value = 30126
$ws.Cells.Item(22, 3).Value = 66
$ws.Cells.Item(22, 4).Value = 25700
$ws.Cells.Item(22, 5).Value = 2616
$ws.Cells.Item(22, 6).Value = 121
$ws.Cells.Item(22, 7).Value = 5
$ws.Cells.Item(22, 8).Value = 1810

$ws.Cells.Item(28, 2).Value = 22385
$ws.Cells.Item(28, 3).Value = 137
$ws.Cells.Item(28, 5).Value = 3872
$ws.Cells.Item(28, 7).Value = 28
$ws.Cells.Item(28, 8).Value = 1403

$ws.Cells.Item(60, 2).Value = 4578
$ws.Cells.Item(60, 3).Value = 156
$ws.Cells.Item(60, 5).Value = 3030

$ws.Cells.Item(71, 2).Value = 2543
$ws.Cells.Item(71, 3).Value = 63
$ws.Cells.Item(71, 4).Value = 1626
$ws.Cells.Item(71, 5).Value = 815

$ws.Cells.Item(98, 4).Value = 400
$ws.Cells.Item(98, 5).Value = 474

$ws.Cells.Item(101, 2).Value = 805
$ws.Cells.Item(101, 3).Value = 8
$ws.Cells.Item(101, 5).Value = 564

$ws.Cells.Item(122, 1).Value = "Paraguay"
$ws.Cells.Item(122, 2).Value = 462
$ws.Cells.Item(122, 3).Value = 22
$ws.Cells.Item(122, 4).Value = 148
$ws.Cells.Item(122, 5).Value = 304
$ws.Cells.Item(122, 6).Value = 9
$ws.Cells.Item(122, 7).Value = 0
$ws.Cells.Item(122, 8).Value = 10

$ws.Cells.Item(123, 1).Value = "Tayikistan"
$ws.Cells.Item(123, 2).Value = 461
$ws.Cells.Item(123, 3).Value = 82
$ws.Cells.Item(123, 4).Value = 0
$ws.Cells.Item(123, 5).Value = 449
$ws.Cells.Item(123, 6).Value = 0
$ws.Cells.Item(123, 7).Value = 4
$ws.Cells.Item(123, 8).Value = 12

$ws.Cells.Item(138, 1).Value = "Madagascar"
$ws.Cells.Item(138, 2).Value = 193
$ws.Cells.Item(138, 3).Value = 35
$ws.Cells.Item(138, 4).Value = 101
$ws.Cells.Item(138, 5).Value = 92
$ws.Cells.Item(138, 8).Value = 0

$ws.Cells.Item(139, 1).Value = "Etiopia"
$ws.Cells.Item(139, 2).Value = 191
$ws.Cells.Item(139, 3).Value = 29
$ws.Cells.Item(139, 4).Value = 93
$ws.Cells.Item(139, 5).Value = 94
$ws.Cells.Item(139, 6).Value = 1
$ws.Cells.Item(139, 8).Value = 4

$ws.Cells.Item(140, 1).Value = "Islas Feroe"
$ws.Cells.Item(140, 2).Value = 187
$ws.Cells.Item(140, 4).Value = 185
$ws.Cells.Item(140, 5).Value = 2
$ws.Cells.Item(140, 6).Value = 0
$ws.Cells.Item(140, 8).Value = 0

$ws.Cells.Item(141, 1).Value = "Martinica"
$ws.Cells.Item(141, 2).Value = 182
$ws.Cells.Item(141, 4).Value = 83
$ws.Cells.Item(141, 5).Value = 85
$ws.Cells.Item(141, 6).Value = 5
$ws.Cells.Item(141, 8).Value = 14

$ws.Cells.Item(142, 1).Value = "Liberia"
$ws.Cells.Item(142, 2).Value = 178
$ws.Cells.Item(142, 3).Value = 0
$ws.Cells.Item(142, 4).Value = 75
$ws.Cells.Item(142, 5).Value = 83
$ws.Cells.Item(142, 8).Value = 20

$ws.Cells.Item(143, 1).Value = "Birmania"
$ws.Cells.Item(143, 2).Value = 176
$ws.Cells.Item(143, 3).Value = 15
$ws.Cells.Item(143, 4).Value = 62
$ws.Cells.Item(143, 5).Value = 108
$ws.Cells.Item(143, 8).Value = 6

$ws.Cells.Item(144, 1).Value = "Santo Tome y Principe"
$ws.Cells.Item(144, 2).Value = 174
$ws.Cells.Item(144, 4).Value = 4
$ws.Cells.Item(144, 5).Value = 167
$ws.Cells.Item(144, 8).Value = 3

$ws.Cells.Item(145, 1).Value = "Republica del Chad"
$ws.Cells.Item(145, 2).Value = 170
$ws.Cells.Item(145, 4).Value = 43
$ws.Cells.Item(145, 5).Value = 110
$ws.Cells.Item(145, 6).Value = 0
$ws.Cells.Item(145, 8).Value = 17

$ws.Cells.Item(146, 1).Value = "Suazilandia"
$ws.Cells.Item(146, 3).Value = 30
$ws.Cells.Item(146, 4).Value = 12
$ws.Cells.Item(146, 5).Value = 139
$ws.Cells.Item(146, 6).Value = 0
$ws.Cells.Item(146, 8).Value = 2

$ws.Cells.Item(147, 1).Value = "Zambia"
$ws.Cells.Item(147, 2).Value = 153
$ws.Cells.Item(147, 3).Value = 7
$ws.Cells.Item(147, 4).Value = 103
$ws.Cells.Item(147, 5).Value = 46
$ws.Cells.Item(147, 6).Value = 1
$ws.Cells.Item(147, 8).Value = 4

$ws.Cells.Item(148, 1).Value = "Guadalupe"
$ws.Cells.Item(148, 2).Value = 152
$ws.Cells.Item(148, 4).Value = 104
$ws.Cells.Item(148, 5).Value = 35
$ws.Cells.Item(148, 6).Value = 4
$ws.Cells.Item(148, 8).Value = 13

$ws.Cells.Item(149, 1).Value = "Gibraltar"
$ws.Cells.Item(149, 2).Value = 144
$ws.Cells.Item(149, 3).Value = 0
$ws.Cells.Item(149, 4).Value = 136
$ws.Cells.Item(149, 5).Value = 8
$ws.Cells.Item(149, 6).Value = 0
$ws.Cells.Item(149, 8).Value = 0

$ws.Cells.Item(150, 1).Value = "Brunei"
$ws.Cells.Item(150, 2).Value = 141
$ws.Cells.Item(150, 3).Value = 2
$ws.Cells.Item(150, 4).Value = 131
$ws.Cells.Item(150, 5).Value = 9
$ws.Cells.Item(150, 6).Value = 2
$ws.Cells.Item(150, 8).Value = 1

$ws.Cells.Item(151, 1).Value = "Benin"
$ws.Cells.Item(151, 2).Value = 140
$ws.Cells.Item(151, 3).Value = 44
$ws.Cells.Item(151, 4).Value = 53
$ws.Cells.Item(151, 5).Value = 85
$ws.Cells.Item(151, 8).Value = 2

$ws.Cells.Item(152, 1).Value = "Guayana Francesa"
$ws.Cells.Item(152, 2).Value = 138
$ws.Cells.Item(152, 4).Value = 112
$ws.Cells.Item(152, 5).Value = 25
$ws.Cells.Item(152, 8).Value = 1

$ws.Cells.Item(153, 1).Value = "Togo"
$ws.Cells.Item(153, 2).Value = 128
$ws.Cells.Item(153, 4).Value = 77
$ws.Cells.Item(153, 5).Value = 42
$ws.Cells.Item(153, 8).Value = 9

$ws.Cells.Item(159, 5).Value = 9
$ws.Cells.Item(159, 7).Value = 1
$ws.Cells.Item(159, 8).Value = 3

$ws.Cells.Item(165, 1).Value = "Liechtenstein"
$ws.Cells.Item(165, 4).Value = 55
$ws.Cells.Item(165, 5).Value = 26
$ws.Cells.Item(165, 6).Value = 0
$ws.Cells.Item(165, 8).Value = 1

$ws.Cells.Item(166, 1).Value = "Barbados"
$ws.Cells.Item(166, 4).Value = 53
$ws.Cells.Item(166, 5).Value = 22
$ws.Cells.Item(166, 6).Value = 4
$ws.Cells.Item(166, 8).Value = 7

$ws.Cells.Item(205, 1).Value = "Montserrat"
$ws.Cells.Item(205, 4).Value = 7
$ws.Cells.Item(205, 6).Value = 1
$ws.Cells.Item(205, 8).Value = 1

$ws.Cells.Item(206, 1).Value = "Seychelles"
$ws.Cells.Item(206, 4).Value = 8
$ws.Cells.Item(206, 6).Value = 0
$ws.Cells.Item(206, 8).Value = 0
